# "Generate Report for Handback" -- refresh the localization-status report
# after a handback: flip each language's Status from "Ready for handoff" to
# "Handed back: in sync with en-US", stamp fresh handback datetimes, clear
# the now-stale "handback file is not the latest" error details, and widen
# a couple of report columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns for both rows ---
$overview.Range("E2:F3").Value = $newStatus

# --- zh-cn detail sheet ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K2").Value = "2016-10-25 02:23:03"
$zhcn.Range("K3").Value = "2016-10-25 02:23:03"
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

# --- de-de detail sheet ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
$dede.Range("K2").Value = "2016-10-25 02:23:20"
$dede.Range("K3").Value = "2016-10-25 02:23:20"
$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

# --- Column widths: widen the Status-bearing columns, shrink Error Detail ---
$overview.Columns.Item(5).ColumnWidth = 29.085322
$overview.Columns.Item(6).ColumnWidth = 29.085322

$zhcn.Columns.Item(3).ColumnWidth = 29.085322
$zhcn.Columns.Item(16).ColumnWidth = 12.750584

$dede.Columns.Item(3).ColumnWidth = 29.085322
$dede.Columns.Item(16).ColumnWidth = 12.750584
